{"js": "// BV_VIP-88: changed from shared_cmd_idx to get_last_received_cmd_idx() in VVC QRs\n//\n// The quick-reference example snippet reads:\n//     v_cmd_idx := shared_cmd_idx;\n// and must become:\n//     v_cmd_idx := get_last_received_cmd_idx(AXISTREAN_VVCT, 1);\n\nconst target = context.document.body.search(\"shared_cmd_idx\", { matchCase: true });\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length === 0) {\n  throw new Error(\"Could not find 'shared_cmd_idx' in the document body.\");\n}\n\n// Replace only the identifier itself; the surrounding ' := ' / ';' text\n// (already split across separate runs) is left exactly as-is, so the\n// resulting visible line reads:\n//   v_cmd_idx := get_last_received_cmd_idx(AXISTREAN_VVCT, 1);\ntarget.items[0].insertText(\"get_last_received_cmd_idx(AXISTREAN_VVCT, 1)\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# BV_VIP-88: changed from shared_cmd_idx to get_last_received_cmd_idx() in VVC QRs\n#\n# The quick-reference example snippet reads:\n#     v_cmd_idx := shared_cmd_idx;\n# and must become:\n#     v_cmd_idx := get_last_received_cmd_idx(AXISTREAN_VVCT, 1);\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"shared_cmd_idx\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$found = $rng.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not find 'shared_cmd_idx' in the document.\"\n}\n\n# Replace only the identifier itself; the surrounding ' := ' / ';' text\n# is left exactly as-is, so the resulting visible line reads:\n#   v_cmd_idx := get_last_received_cmd_idx(AXISTREAN_VVCT, 1);\n$rng.Text = \"get_last_received_cmd_idx(AXISTREAN_VVCT, 1)\"\n"}
